$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0.2077624555612785"
$ws.Range("C2").Value = [double]"1.704538857752482"

$ws.Range("B3").Value = [double]"9.908539415492668e-05"
$ws.Range("C3").Value = [double]"7.983629470097142"

$ws.Range("B4").Value = [double]"0.3599305001028352"
$ws.Range("C4").Value = [double]"0.9790362244201276"

$ws.Range("B5").Value = [double]"5.955649589091438e-112"
$ws.Range("C5").Value = [double]"1.704538857752482"
